$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell C2 ("Vaibhav@321", hyperlinked) is being removed entirely.
# The remaining hyperlinked cells (B2, C3, B4, C4, C5) keep the same
# targets/labels as before; only C2's hyperlink relationship disappears.
#
# This runtime only supports removing hyperlinks in bulk (Worksheet.Hyperlinks.Delete),
# so capture the formatting of the surviving hyperlinked cells first, wipe every
# hyperlink, clear C2's content, then restore the 5 hyperlinks that remain along
# with the original cell formatting.

$keep = @(
    @{ Addr = "B2"; Target = "mailto:v.malkar@smartshiphub.com"; Text = "v.malkar@smartshiphub.com" },
    @{ Addr = "C3"; Target = "mailto:Vaibhav@321";               Text = "Vaibhav@321" },
    @{ Addr = "B4"; Target = "mailto:v.malkar@smartshiphub.com"; Text = "v.malkar@smartshiphub.com" },
    @{ Addr = "C4"; Target = "mailto:Vaibhav@3211";              Text = "Vaibhav@3211" },
    @{ Addr = "C5"; Target = "mailto:Vaibhav@3211";              Text = "Vaibhav@3211" }
)

# NOTE: mutating hashtable entries through a `foreach` loop variable does not
# persist back into the backing array on this host - use indexed `for` loops
# so the captured formatting survives between the two passes below.
for ($i = 0; $i -lt $keep.Count; $i++) {
    $f = $ws.Range($keep[$i].Addr).Font
    $keep[$i].FontName = $f.Name
    $keep[$i].FontSize = $f.Size
    $keep[$i].FontUnderline = $f.Underline
}

# Remove every hyperlink on the sheet (bulk operation is the only reliable one here).
$ws.Hyperlinks.Delete()

# Drop the value that lived in C2 along with its (now-gone) hyperlink.
$ws.Range("C2").ClearContents()

# Re-create the hyperlinks that should still exist, preserving document order
# so relationship ids come out as rId1..rId5 in the same B2,C3,B4,C4,C5 order.
for ($i = 0; $i -lt $keep.Count; $i++) {
    $ws.Hyperlinks.Add($ws.Range($keep[$i].Addr), $keep[$i].Target, [Type]::Missing, [Type]::Missing, $keep[$i].Text)

    $f = $ws.Range($keep[$i].Addr).Font
    $f.Name = $keep[$i].FontName
    $f.Size = $keep[$i].FontSize
    $f.Underline = $keep[$i].FontUnderline
    $f.ThemeColor = 1
    $f.TintAndShade = 0
}

# Update the active selection to E13
$ws.Range("E13").Select()
